$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.866.82"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.832.67"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'245.07"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'0.6903"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.07707"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "'0.3048"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").Value = "'23.36"
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("D11").Value = "'0.07809"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "1.833.78"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "'5.091"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "'91.15"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "'0.6814"
$ws.Range("E15").Value = "  -2.62%  "
$ws.Range("D16").Value = "'6.418"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "'0.000008309"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D18").Value = "28.881.50"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "'242.02"
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").Value = "2.075.15"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").Value = "'12.71"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'7.451"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'0.1478"
$ws.Range("D26").Value = "'158.31"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").Value = "'8.791"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "'18.23"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "'4.220"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").Value = "'4.151"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").Value = "'1.194"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").Value = "'0.05101"
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("D34").Value = "'0.7805"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("D35").Value = "'1.853"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "'2.689"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D39").Value = "1.224.44"
$ws.Range("E39").Value = "  -3.93%  "
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").Value = "'0.9573"
$ws.Range("E41").Value = "  +6.47%  "
$ws.Range("D42").Value = "'109.10"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "'5.856"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'9.614"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("D47").Value = "1.977.08"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").Value = "'0.5158"
$ws.Range("D49").Value = "'64.01"
$ws.Range("E49").Value = "  -9.28%  "
$ws.Range("D50").Value = "'1.749"
$ws.Range("E50").Value = "  -2.69%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'6.920"
$ws.Range("E51").Value = "  -1.78%  "
